# Applies the updated crypto market snapshot (prices / 1h volume %) to sheet1.
# Values are written with a leading apostrophe to force text storage (matching
# the source workbook, where these columns are inline strings, not numbers),
# then the style is reset to "Normal" so no stray number-format / quote-prefix
# style gets attached to the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'37.188.39"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +1.35%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'2.006.60"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +2.16%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  +0.12%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'246.76"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +1.03%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'0.623"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +1.21%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'59.88"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  -1.74%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = "'  +0.04%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Value = "'  +2.14%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.0811"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  +2.04%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("E11").Value = "'  +0.50%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'15.05"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  +5.65%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'22.34"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +2.18%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'2.301.29"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +2.21%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("E15").Value = "'  +0.08%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("E16").Value = "'  +2.64%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'1.999.18"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +1.82%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'37.081.66"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +1.20%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'70.24"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +0.37%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'0.0₃0866"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  +1.29%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'5.20"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +1.71%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("E22").Value = "'  -0.03%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("E23").Value = "'  -0.07%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("E24").Value = "'  -0.09%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("E25").Value = "'  +0.37%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'9.44"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  +1.92%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'164.64"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  +2.28%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'0.137"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -5.37%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'19.67"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  +0.88%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("E30").Value = "'  +14.07%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("E31").Value = "'  +1.12%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("E32").Value = "'  +0.51%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'0.0656"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  +5.98%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'4.47"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  +0.70%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'2.44"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  +6.77%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("E36").Value = "'  +0.08%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("B37").Value = "'RenderToken"
$ws.Range("B37").Style = "Normal"
$ws.Range("C37").Value = "'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("C37").Style = "Normal"
$ws.Range("D37").Value = "'3.40"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  -5.44%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("B38").Value = "'WEMIXToken"
$ws.Range("B38").Style = "Normal"
$ws.Range("C38").Value = "'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("C38").Style = "Normal"
$ws.Range("D38").Value = "'1.81"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  +2.67%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'5.36"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -4.80%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.0984"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -0.01%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("E41").Value = "'  +0.90%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("E42").Value = "'  +0.77%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("E43").Value = "'  +1.10%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'16.63"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +2.37%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'91.78"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +3.52%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'1.371.17"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +0.07%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Value = "'  +1.09%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'7.38"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  +2.81%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'2.09"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  +14.01%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'46.90"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +5.53%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("E51").Value = "'  -0.34%  "
$ws.Range("E51").Style = "Normal"
